# Landscaping Data.xlsx - "Add files via upload" edit
#
# Appends 14 new observation rows (rows 366-379) to Sheet1, continuing the
# existing data table (which previously ended at row 365). Columns:
#   A Date | B Plant_Type | C Plant_Size | D Low | E High | F Temp_Diff (=ABS(D-E))
#   G Rain | H Growth | I Pruned | J Quadrant | K Shade | L UV | M Humidity
#   N Dew_Point | O Pressure | P Wind_Gust | Q Cloud_Cover | R Visibility
#   S AQI | T Pollen

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows' data (column F, Temp_Diff, is a formula and is filled in separately).
$data = @(
    @(45839, "Flowering", "Large", 67, 78, 0.5, 0.2, "No", 2, "Dark", 6, 0.67, 70, 29.99, 13, 0.78, 8.6999999999999993, 49, 0),
    @(45839, "Nonflowering", "Medium", 67, 78, 0.5, 0.1, "No", 3, "Dark", 6, 0.67, 70, 29.99, 13, 0.78, 8.6999999999999993, 49, 0),
    @(45839, "Nonflowering", "Small", 67, 78, 0.5, 0.2, "No", 3, "Neutral", 6, 0.67, 70, 29.99, 13, 0.78, 8.6999999999999993, 49, 0),
    @(45839, "Nonflowering", "Medium", 67, 78, 0.5, 0.25, "No", 3, "Neutral", 6, 0.67, 70, 29.99, 13, 0.78, 8.6999999999999993, 49, 0),
    @(45839, "Nonflowering", "Medium", 67, 78, 0.5, 0.3, "No", 3, "Bright", 6, 0.67, 70, 29.99, 13, 0.78, 8.6999999999999993, 49, 0),
    @(45839, "Nonflowering", "Large", 67, 78, 0.5, 0.3, "No", 4, "Bright", 6, 0.67, 70, 29.99, 13, 0.78, 8.6999999999999993, 49, 0),
    @(45839, "Tree", "Medium", 67, 78, 0.5, 1.1000000000000001, "No", 1, "Neutral", 6, 0.67, 70, 29.99, 13, 0.78, 8.6999999999999993, 49, 0),
    @(45840, "Flowering", "Large", 65, 84, 0, 0, "No", 2, "Neutral", 9, 0.57999999999999996, 66, 29.99, 9, 0.36, 9.9, 45, 0),
    @(45840, "Nonflowering", "Medium", 65, 84, 0, 0, "No", 3, "Neutral", 9, 0.57999999999999996, 66, 29.99, 9, 0.36, 9.9, 45, 0),
    @(45840, "Nonflowering", "Small", 65, 84, 0, 0.1, "No", 3, "Bright", 9, 0.57999999999999996, 66, 29.99, 9, 0.36, 9.9, 45, 0),
    @(45840, "Nonflowering", "Medium", 65, 84, 0, 0.1, "No", 3, "Bright", 9, 0.57999999999999996, 66, 29.99, 9, 0.36, 9.9, 45, 0),
    @(45840, "Nonflowering", "Medium", 65, 84, 0, 0.25, "No", 3, "Bright", 9, 0.57999999999999996, 66, 29.99, 9, 0.36, 9.9, 45, 0),
    @(45840, "Nonflowering", "Large", 65, 84, 0, 0.2, "No", 4, "Bright", 9, 0.57999999999999996, 66, 29.99, 9, 0.36, 9.9, 45, 0),
    @(45840, "Tree", "Medium", 65, 84, 0, 0.5, "No", 1, "Neutral", 9, 0.57999999999999996, 66, 29.99, 9, 0.36, 9.9, 45, 0)
)

$firstNewRow = 366
$lastExistingRow = 365
$lastNewRow = $firstNewRow + $data.Count - 1

# Carry the existing row's formatting (incl. the m/d/yyyy date style on col A)
# down across the whole new block before writing values.
$ws.Range("A" + $lastExistingRow + ":T" + $lastExistingRow).Copy()
$ws.Range("A" + $firstNewRow + ":T" + $lastNewRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $firstNewRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]        # A Date
    $ws.Cells.Item($r, 2).Value = $row[1]         # B Plant_Type
    $ws.Cells.Item($r, 3).Value = $row[2]         # C Plant_Size
    $ws.Cells.Item($r, 4).Value = $row[3]         # D Low
    $ws.Cells.Item($r, 5).Value = $row[4]         # E High
    $ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"   # F Temp_Diff
    $ws.Cells.Item($r, 7).Value = $row[5]         # G Rain
    $ws.Cells.Item($r, 8).Value = $row[6]         # H Growth
    $ws.Cells.Item($r, 9).Value = $row[7]         # I Pruned
    $ws.Cells.Item($r, 10).Value = $row[8]        # J Quadrant
    $ws.Cells.Item($r, 11).Value = $row[9]        # K Shade
    $ws.Cells.Item($r, 12).Value = $row[10]       # L UV
    $ws.Cells.Item($r, 13).Value = $row[11]       # M Humidity
    $ws.Cells.Item($r, 14).Value = $row[12]       # N Dew_Point
    $ws.Cells.Item($r, 15).Value = $row[13]       # O Pressure
    $ws.Cells.Item($r, 16).Value = $row[14]       # P Wind_Gust
    $ws.Cells.Item($r, 17).Value = $row[15]       # Q Cloud_Cover
    $ws.Cells.Item($r, 18).Value = $row[16]       # R Visibility
    $ws.Cells.Item($r, 19).Value = $row[17]       # S AQI
    $ws.Cells.Item($r, 20).Value = $row[18]       # T Pollen
}

# Match the author's final on-screen selection/scroll state.
$ws.Range("Q373:Q379").Select()
$excel.ActiveWindow.ScrollRow = 355
$excel.ActiveWindow.ScrollColumn = 1
